$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows (23-26) need the same header/data style ("s=1") applied to column A
# as the rest of the table before setting their values.
$ws.Range("A2").Copy()
$ws.Range("A23:A26").PasteSpecial(-4122)

$names = @("model_18_5_0","model_18_5_22","model_18_5_21","model_18_5_20","model_18_5_19","model_18_5_18","model_18_5_17","model_18_5_16","model_18_5_15","model_18_5_14","model_18_5_13","model_18_5_23","model_18_5_12","model_18_5_10","model_18_5_9","model_18_5_8","model_18_5_7","model_18_5_6","model_18_5_5","model_18_5_4","model_18_5_3","model_18_5_2","model_18_5_1","model_18_5_11","model_18_5_24")
$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q")
$vals = @(0.9298247212236502,0.7246052329558914,0.9778201568565887,0.8097653136555514,0.8648662959483421,0.469262412304898,1.841566075300257,0.0177496122277409,0.9719354387791824,0.4948424447787041,0.3376508542801637,0.6850273077074358,1.129554360817876,0.7141903383985572,75.51318630499273,120.6115918251162)

for ($i = 0; $i -lt $names.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $names[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Range($cols[$j] + $r).Value = $vals[$j]
    }
}
